$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(20, 1).Value = '(Clara) Catel cu telecomanda'
$ws.Cells.Item(20, 2).Value = 'https://m.media-amazon.com/images/I/71A+tWCmVlL._AC_SX679_.jpg'
$ws.Cells.Item(20, 3).Value = 'https://www.amazon.de/-/en/furReal-Walkin-Interactive-Bouncing-Effects/dp/B0CRBJ8JFR/ref=sr_1_2?crid=1VQJRA7B0MIQV&dib=eyJ2IjoiMSJ9.53Uu6FxJ4b6G8ZdU8vF8f5_69fBrUXbiB0mmwb1j4o3FN5oTVqD_WxOQ-IU9ZCdLyRupmmDs2Bni1O65n56dwqXUKAQ_hfty5U3xd6uubP6GEKV_XpgF4J5s--YdZWRl4b9kz0DjalmySz_sZaxBs7TKoEQrpKoOIm0HBvMOyJci82LMGrfxTtnRQVckEemvsf_-Ukh732hHa8ZrScNKi2kpUIJZSlMODf42dJNSYZY.WOEefNOMKp6eIVjyUoZ4k5cSLQjkeM2JF-DkR-rmPyE&dib_tag=se&keywords=furreal+friends+hund+gogo&qid=1729516311&sprefix=furreal+gogo+hun%2Caps%2C104&sr=8-2'
$ws.Cells.Item(20, 4).Value = '37 EUR'

$ws.Cells.Item(21, 1).Value = '(Clara) Telefon Gabby'
$ws.Cells.Item(21, 2).Value = 'https://m.media-amazon.com/images/I/71L5GUCfx4L._AC_SX679_.jpg'
$ws.Cells.Item(21, 3).Value = 'https://www.amazon.de/-/en/VTech-Gabbys-Dollhouse-Learning-Mobile/dp/B0CKRQ3N9W/ref=sr_1_2?crid=1TGUUBO2GUPXG&dib=eyJ2IjoiMSJ9.y1-F9U8wTKploOwLbra7x_dsi6XA_61WOmjL6YdtfKrMZd5PgbtYetS_xvrxVlN5P2N30MVC6N3-x3f18WvJLvsACyrQPzED_H1v5AXbB5G5syHB0uEk7crueAZVtpCq1DlsW_jJNYq258-TK7XttLHSAv_SDIynj2wHjEpJ8M9MAHFKYdOwMS1MVOfpFPOYGVf7604XryxZFTUkz68xV3oE47bLD-NFwuB1HEwlYbg.TibLhJEUOcQ47hE3d_W_l7dd4y7OSJytlm989H5OdP4&dib_tag=se&keywords=gabby+lernhandy&qid=1729517077&sprefix=gabby+lern%2Caps%2C108&sr=8-2#customerReviews'
$ws.Cells.Item(21, 4).Value = '20 EUR'

$ws.Cells.Item(22, 1).Value = '(Clara) Plus Stitch cu sunet si miros'
$ws.Cells.Item(22, 2).Value = 'https://m.media-amazon.com/images/I/81wmHRHoM8L._AC_SX679_.jpg'
$ws.Cells.Item(22, 3).Value = 'https://www.amazon.de/-/en/Stitch-Disney-Sound-Scent-Small/dp/B0C6YL2H1D/ref=sr_1_1?crid=6EMXNMGZ24N&dib=eyJ2IjoiMSJ9.GuA80LmF0yuVPEU-kqsv9uLayaq9ADpBvReu1nPc2asAtEJmgEUUTq10j-altcYb3t9k1YyGDPJevWY5WAl7q5HRvhO8xbA4PHlUbBc_VRiNiXsoIrryF9CCiw_8Ymsuuje7-xai4s4kYksI-f3wNk23ZAm99jSMVpGYeEs9Du2EbsLhExjtRQruu-kR9luH8MB5o7ECPnd7ASgVz6ZW3uqoxwIynSUbp8zjBukwaj0GZ0w34C_HnUySu_HAhnPT-3Wp0wgwUcFnsx4PPpEmbs-Hw7QzcljVB4diYLvaOs0.0Y9HuQtkgUCAVRkvEzI6CmJ_iNtgsbo1m-JlnVZ6QXk&dib_tag=se&keywords=stitch+soundfunktion+duft+15cm&qid=1729516484&sprefix=stitch+soundfunktion+duft+15cm%2Caps%2C76&sr=8-1'
$ws.Cells.Item(22, 4).Value = '17 EUR'

$ws.Cells.Item(23, 1).Value = '(Clara) Set potiuni'
$ws.Cells.Item(23, 2).Value = 'https://m.media-amazon.com/images/I/81yFvCyxfLL._AC_SX679_.jpg'
$ws.Cells.Item(23, 3).Value = 'https://www.amazon.de/-/en/Mixies-Cauldron-Purple-Interactive-Conjure/dp/B0BX6YW5N9/ref=sr_1_1?crid=358IKESIBXVKO&dib=eyJ2IjoiMSJ9.oFY88I1D4MFo5A9ePZtzXjbQS4l-SpFJkwgAbWg7MMhtOcENlN6wulpw416W8xsS2fLUPGO47TnycXTPWIxHE2CDymxryzuQt0W7vBohROKY7k9Gxd5wxAniONIj00EQFSdEFqwA14QjSRc8QUbA489aooXS7DEeynFt9EmUOEY8HuoiqkTF2L2tH2U6weSLdM-zRruSR3jTRicqwIRaJ7KS33hoqDz6YFgspU9AJP0G5bDKdrcfvCTO2WeSnZDhWVx8KzBe_NmTT37oDFQSA7zNCmpXIAokJ02pSwIx80U.YJEOHI_2v0L6eI_rOrgvMkf9kcId4bmGS2Bwd1uJk7k&dib_tag=se&keywords=my+magic+mixies&qid=1729516398&sprefix=my+magic+mix%2Caps%2C99&sr=8-1'
$ws.Cells.Item(23, 4).Value = '66 EUR'

$ws.Cells.Item(24, 1).Value = 'Gravitrax Junior Starter Set L - Ice'
$ws.Cells.Item(24, 2).Value = 'https://m.media-amazon.com/images/I/715c7oaL5bL._AC_SX679_.jpg'
$ws.Cells.Item(24, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-GraviTrax-Junior-Starter-Environmentally/dp/B0BSX847R3/ref=sr_1_2?crid=1X7LU8OFXMUDY&dib=eyJ2IjoiMSJ9.HacZPOrg5jUAw8zb8YX2rP9wurG1VPUH9ajjsfsttYOMungu-BuQav_M7W7LDs_HaLJiA4CCD1XsLUL7P628PHiaWR93cn8kXvmZKpxugWklhxzMUdAxO8NonF9N-uGQB4B7-4lX_3SD91JPzXuOLiPvsGDu3Menod9Pg6gkgWDHxOuh2ZRNrm5yFve3aLzysRcmg5diMQ92Yws4UuWjJvEXvD77q8MLTFxrgWAXOErjvtLv-E_gwG34fWMeduft6MSmf5F-oh8goqjjolGhi-h-dQQ5YecsoFCB3ncrJys.9EKnkXvyxCPIKZXs9Mx9U7m2t-LMTsbMy9fr4BxJK6M&dib_tag=se&keywords=gravitrax+junior+starterset&qid=1729516745&sprefix=gravitrax+junior%2Caps%2C98&sr=8-2'
$ws.Cells.Item(24, 4).Value = '36 EUR'

$ws.Cells.Item(25, 1).Value = 'Gravitrax Junior Starter Set Frozen'
$ws.Cells.Item(25, 2).Value = 'https://m.media-amazon.com/images/I/81yqCdPH3SL._AC_SX679_.jpg'
$ws.Cells.Item(25, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-73385-GraviTrax-Starter-Children/dp/B0CSZ63YV3/ref=sr_1_5?crid=1X7LU8OFXMUDY&dib=eyJ2IjoiMSJ9.HacZPOrg5jUAw8zb8YX2rP9wurG1VPUH9ajjsfsttYOMungu-BuQav_M7W7LDs_HaLJiA4CCD1XsLUL7P628PHiaWR93cn8kXvmZKpxugWklhxzMUdAxO8NonF9N-uGQB4B7-4lX_3SD91JPzXuOLiPvsGDu3Menod9Pg6gkgWDHxOuh2ZRNrm5yFve3aLzysRcmg5diMQ92Yws4UuWjJvEXvD77q8MLTFxrgWAXOErjvtLv-E_gwG34fWMeduft6MSmf5F-oh8goqjjolGhi-h-dQQ5YecsoFCB3ncrJys.9EKnkXvyxCPIKZXs9Mx9U7m2t-LMTsbMy9fr4BxJK6M&dib_tag=se&keywords=gravitrax+junior+starterset&qid=1729516745&sprefix=gravitrax+junior%2Caps%2C98&sr=8-5'
$ws.Cells.Item(25, 4).Value = '36 EUR'

$ws.Cells.Item(26, 1).Value = '(Clara) Papusa '
$ws.Cells.Item(26, 2).Value = 'https://m.media-amazon.com/images/I/51WGSlK3v6L._AC_SX679_.jpg'
$ws.Cells.Item(26, 3).Value = 'https://www.amazon.de/-/en/105733519-Princess-Movable-Suitable-Children/dp/B08QYTD24X/ref=sr_1_2?crid=K5EE4ZX6PC2F&dib=eyJ2IjoiMSJ9.6WxTQND8ZRjVlveNJY2qRw0Mm4LH90lyvKyzz7nDafVgwRCnIYzN0bWBiMWZiMP8P163-QVHZGnO9vDq3RJAZwEAzmDQXfBGsgqWPa4-0g7KruRsNtj8nfuKPUBWSzuOZ-3hLnN7hABHEDEYMhLbj_OGraTMT0BRQWcsBS1khICmFnPU3eHBo41bmEBw0YXRe8BbNMYczPvSacihlWG3O1IyUVgp8i5f3cKO4s39KRVhD4gPPFJEUYL93RC4pV6NOYax081G7OZyfV0wcViHLvD-wO1o6ghtjUnqc7IxnQk.KN4nYMyUPUc1xxm5YUWrMV7Lym367vaok8RY7lT9OrI&dib_tag=se&keywords=steffi+love+pferd&qid=1729591110&sprefix=stefi+lo%2Caps%2C116&sr=8-2#customerReviews'
$ws.Cells.Item(26, 4).Value = '23 EUR'

$ws.Cells.Item(27, 1).Value = '(Maia) Plus Angel 25cm'
$ws.Cells.Item(27, 2).Value = 'https://m.media-amazon.com/images/I/71d-QML9QWL._SX450_.jpg'
$ws.Cells.Item(27, 3).Value = 'https://www.amazon.de/-/en/6315876954NPB-Disney-Stitch-Figure-Suitable/dp/B07Q271JXM/ref=sr_1_5?crid=241IRMM0W9185&dib=eyJ2IjoiMSJ9.1CdE7TiRnOOhoEpwc0JJyS_0nYoXwSwrMblaLAz0NGE7ndxcZSKt8RfhJ3YxOyHSuP2BuuGvwzCL5rkgkGWlTEGwrYA6Iunn1kbqmIo4kIfzPyMHq8PKN0B_NvsiZgJMU6Cq1CgfZvVE1slDBKy55ac2JdUC0xtYL7TrVEAZppMm0_iOrTwRJAhWCBuhzyg2_BiqAUfRoUsEry4QSv1IUJ050D-Y4pXqMkjEmYaV-IWh3WGbZmYhjIzyLvVKyF1_gUFBJmxamD0RowRFFwB8XMg4gdY0fJ4ycsLJ2WpphmQ.lxXY7_J0pJABvnJCYEqHGr9RG7wKTz8lgYLLIVqndFE&dib_tag=se&keywords=disney%2Bstitch%2Bleroy%2Bangel%2B25cm&qid=1729845349&s=toys&sprefix=disney%2Bstitch%2Bleroy%2Bangel%2B25cm%2Ctoys%2C107&sr=1-5&th=1'
$ws.Cells.Item(27, 4).Value = '17 EUR'

$ws.Cells.Item(28, 1).Value = '(Maia) Telefon Paw Patrol'
$ws.Cells.Item(28, 2).Value = 'https://m.media-amazon.com/images/I/81FXrd8-pOL._AC_SY879_.jpg'
$ws.Cells.Item(28, 3).Value = 'https://www.amazon.de/-/en/VTech-Paw-Patrol-Learning-Educational/dp/B09W2RPFGM/ref=sr_1_1?crid=27879TVAXSDMC&dib=eyJ2IjoiMSJ9.gBFhkTT4VJsbXzfmFH6nLA9XHBnrSXZSr15G23ofAmSvatdZ4x17fOGRi7w4pfNnI8Dtq03vrJ9vfalk2pfMjrR2aBn3iR8WlgMtWFQZVF3qN-4rwWSeDt4KtrpfAwp0KSmPhDMPiN7kIrN-Va8ruKlWoAMbw7MUPxuouDJ9AYybPTftjGIHt3nmNsPoqOkKAE6tcHugBe2dJnksPjcak0_V8XY-vgAjwHWEBSS81yMKYcJuo8ID6RWHJ2OpxxxssIweUUU8036BkNdUV7owvjR48cMUbb-CI359LSXdDnI.7IO2Hjo0k8F9hlZkI2mQUWaXtRlBqVRgiaV0klHHsuw&dib_tag=se&keywords=paw+patrol+lernhandy&qid=1729845016&s=toys&sprefix=paw+patrol+lernhandy%2Ctoys%2C103&sr=1-1'
$ws.Cells.Item(28, 4).Value = '25 EUR'

$ws.Cells.Item(29, 1).Value = '(Maia) Sirena Barbie'
$ws.Cells.Item(29, 2).Value = 'https://m.media-amazon.com/images/I/619HvbXSNlL._AC_SX679_.jpg'
$ws.Cells.Item(29, 3).Value = 'https://www.amazon.de/Accessories-Unpacking-Surprises-Mermaid-HRK12/dp/B0CRLYSY7H/ref=sr_1_1_sspa?crid=3DH370C6RA19Q&dib=eyJ2IjoiMSJ9.5n_kquIaPyKlHSzp0JrBDi1QPncmHENqyPIErFgNCekrmsYVJztWt4_iaiyT2lpUl_fQjMaczD4LNk09OHLPLD3kGxzFD3v7VcmMWwZ7y4SK2LEzgnkCNxGLWQ4imixY5oDypJFl34CKCmcVmhEYC5Fg5dDtJg1h3XPtGPIV0dOsAg4Ao43jfrXjwZ1fgLr_t_CmguQRkMJeblSD4JqhoHb2s4gQr1tnw1BOxTfAmULwbKy9PqtVYCr6wkY5fl7yfOB5TKlNfP2vqeF6_8wUdCrMe1dHcfR1PamvrsfnU9s.FVX6P4lyQUEfLLVwlFRoiqjdRgP77RduzLy8VyK81-A&dib_tag=se&keywords=barbie+color+reveal&qid=1729845088&s=toys&sprefix=barbie+color+re%2Ctoys%2C111&sr=1-1-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9hdGY&psc=1'
$ws.Cells.Item(29, 4).Value = '15 EUR'

$ws.Cells.Item(30, 1).Value = '(Maia) Trusa doctor'
$ws.Cells.Item(30, 2).Value = 'https://m.media-amazon.com/images/I/81J3vtHatLL._AC_SX679_.jpg'
$ws.Cells.Item(30, 3).Value = 'https://www.amazon.de/Theo-Klein-10-Piece-Doctors-Case/dp/B0002ML3MK/ref=sr_1_1_sspa?crid=4DLSOMO7CFFH&dib=eyJ2IjoiMSJ9.7nxPsw_LkG-CJICNBTkfGjGYe-WPEO5FoaI3DQrl1kgjO8vkYdeqj0rWaZO9cy94iPuvkJT1NpOZlbZMTuFqDQuvjf_now8n1Ez2X0baqkrWeGDlgHIAA16tcgS1WQHt1JIFSmS9MWy9MnEFjlVrDAQUPRU09L0sNdqwJ9_j3ADt0nGHwtsJD0poIKwlm0hYBZQyL-X2DFiMTDA0jhRW-rjaqbofeXHDlVIPK_Ta4l40Ry9_mfnhiDmKDCfFQUf2FkzNu_H9bnCGO3VpeeSz62xXFrtR9W-gPV19sciwcL4.FD4WElDR3-ulY7A6O0QXKMxLrpVgzMPTc_bBzBppkY4&dib_tag=se&keywords=klein%2Barzt%2Bkoffer&qid=1729845426&sprefix=klein%2Bar%2Caps%2C319&sr=8-1-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9hdGY&th=1'
$ws.Cells.Item(30, 4).Value = '16 EUR'

$ws.Cells.Item(31, 1).Value = '(Maia) Puzzle '
$ws.Cells.Item(31, 2).Value = 'https://m.media-amazon.com/images/I/61CDR8ItfmL._AC_SX679_.jpg'
$ws.Cells.Item(31, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-Childrens-Puzzle-Favourite-Dinosaurs/dp/B09DDB9CHM/ref=sr_1_8?crid=6H9K2UOBYCGW&dib=eyJ2IjoiMSJ9.pOzpXEoCcVDebiGOf1UfBHoQE_98Yt-9ps8MjVv6Q7yBp-VUeLXaNn8tibvb7PTihZDgpwF3CBkgBmMa7_6LHkPj0Bz1-sjw952GkWIKwh3nTvl6_vhQd0uAHC3xo-Ml3mrK_rRdfLVKycfNA-7-yPxtL0w892i6zp_Zb90lI3mhkeavM4x6VR-KpuRXjI7LDzy1DOpnclpT3xF02ync5-oXLmtDXp6PtQGM7PeF8ppXmpgW2PVFyrwnOxRa0MLE0L2C89Rv87NlozWp_QrYiK0SezNgwccy6LURdBcbc6s.pIkkiMw7lwPPFrD9UNIUgO7mHajq5MnbaHl2AdoJJlM&dib_tag=se&keywords=ravensburger+puzzle+ab+3+2x12&qid=1729845613&sprefix=ravensburger+puzzle+ab+3+2x%2Caps%2C97&sr=8-8'
$ws.Cells.Item(31, 4).Value = '8 EUR'

$ws.Cells.Item(32, 1).Value = '(Maia) Puzzle '
$ws.Cells.Item(32, 2).Value = 'https://m.media-amazon.com/images/I/81cnfZUCtTL._AC_SX679_.jpg'
$ws.Cells.Item(32, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-Childrens-Puzzle-Dollhouse-Children/dp/B0BXLLCXSZ/ref=sr_1_20?crid=6H9K2UOBYCGW&dib=eyJ2IjoiMSJ9.pOzpXEoCcVDebiGOf1UfBHoQE_98Yt-9ps8MjVv6Q7yBp-VUeLXaNn8tibvb7PTihZDgpwF3CBkgBmMa7_6LHkPj0Bz1-sjw952GkWIKwh3nTvl6_vhQd0uAHC3xo-Ml3mrK_rRdfLVKycfNA-7-yPxtL0w892i6zp_Zb90lI3mhkeavM4x6VR-KpuRXjI7LDzy1DOpnclpT3xF02ync5-oXLmtDXp6PtQGM7PeF8ppXmpgW2PVFyrwnOxRa0MLE0L2C89Rv87NlozWp_QrYiK0SezNgwccy6LURdBcbc6s.pIkkiMw7lwPPFrD9UNIUgO7mHajq5MnbaHl2AdoJJlM&dib_tag=se&keywords=ravensburger+puzzle+ab+3+2x12&qid=1729845613&sprefix=ravensburger+puzzle+ab+3+2x%2Caps%2C97&sr=8-20'
$ws.Cells.Item(32, 4).Value = '7 EUR'

$ws.Cells.Item(33, 1).Value = '(Maia) Puzzle '
$ws.Cells.Item(33, 2).Value = 'https://m.media-amazon.com/images/I/81fqssNDPmL._AC_SX679_.jpg'
$ws.Cells.Item(33, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-Childrens-Puzzle-Flower-Meadow/dp/B015XS7UAQ/ref=sr_1_6?crid=6H9K2UOBYCGW&dib=eyJ2IjoiMSJ9.pOzpXEoCcVDebiGOf1UfBHoQE_98Yt-9ps8MjVv6Q7yBp-VUeLXaNn8tibvb7PTihZDgpwF3CBkgBmMa7_6LHkPj0Bz1-sjw952GkWIKwh3nTvl6_vhQd0uAHC3xo-Ml3mrK_rRdfLVKycfNA-7-yPxtL0w892i6zp_Zb90lI3mhkeavM4x6VR-KpuRXjI7LDzy1DOpnclpT3xF02ync5-oXLmtDXp6PtQGM7PeF8ppXmpgW2PVFyrwnOxRa0MLE0L2C89Rv87NlozWp_QrYiK0SezNgwccy6LURdBcbc6s.pIkkiMw7lwPPFrD9UNIUgO7mHajq5MnbaHl2AdoJJlM&dib_tag=se&keywords=ravensburger+puzzle+ab+3+2x12&qid=1729845613&sprefix=ravensburger+puzzle+ab+3+2x%2Caps%2C97&sr=8-6'
$ws.Cells.Item(33, 4).Value = '8 EUR'

$ws.Cells.Item(34, 1).Value = 'Lotti Karotti Unicorn'
$ws.Cells.Item(34, 2).Value = 'https://m.media-amazon.com/images/I/81ItnMCR2dL._AC_SX679_.jpg'
$ws.Cells.Item(34, 3).Value = 'https://www.amazon.de/-/en/Ravensburger-22580-Adaptation-Well-known-Entertaining/dp/B0D7W16RHW/ref=sr_1_1?crid=33W3ZJWMCJ3U6&dib=eyJ2IjoiMSJ9.rd6AsEmkfsTPiA2mQ7ZBblS7Q1FWU-AiK9IlqhdQtm1eKJM3npv-dbjOkPmvBWCq9bgO-jw1zOherZVmmaKdaBPHdUL-i_Dnk4ecgsQJ0jaSpZd1vCyN5FswK52_lwM3uFiNk6CxiybFftJXUMU61HC7hVg-h3JZgQBhVNGwh142atAjOX8WID8zdXT0iTHNVowH8jQ1EAfmZ4_7EDnFp2zjTx1Ooe87YiOb-myIG2v1AWLqO2vWDIEYRt-EFCXJiR77cO78moxyrYP9PRUB4MqlM3k2GzBr3ZGjyH7Rvr8.6QMekaHl0LxdB9ggAzZRg8y7E0W2d3cDjQWitdv3phM&dib_tag=se&keywords=lotti+karotti+einhorn&qid=1729845712&sprefix=lotti+%2Caps%2C151&sr=8-1'
$ws.Cells.Item(34, 4).Value = '21 EUR'

$ws.Cells.Item(35, 1).Value = '(Clara) Salon Gabby'
$ws.Cells.Item(35, 2).Value = 'https://m.media-amazon.com/images/I/61tVFCO2ZLL._AC_SX679_.jpg'
$ws.Cells.Item(35, 3).Value = 'https://www.amazon.de/-/en/Smoby-Dressing-Popular-Dollhouse-Gabby/dp/B0CY9K6MRH/ref=sr_1_3?crid=36MMI6PHDMFMC&dib=eyJ2IjoiMSJ9.c9vuC5JfAN0gUngFM1I5YvqRamfXFTIdUWcqhdGS8xiExSwOh1f2Ctj1jGxYtL1xR-7DrJHfI5QYHtYSbIaXQbcIhZJQs9otBz3joQ4rEDXcjOde_OFCE0-c9t34nZxr3iKyP68UAtY_KnlR0ki14sszZkWBM43yS2aYD9QweXBwjeOEuPVwcZubVT9oVUQpDjd_mm3P6LZo9gE2uIP1gE51As9jStDgAgt9NjJAVQH4IJopxCtwawDnx8KVNHfo8WnxL7ZdDpurcveNJhO1r0WVsnhM4AOlzN7lDJyDxGg.T-f_oMEP2Fza0bmztN7nwB0ij_SY87FP4YISnEOcJ2k&dib_tag=se&keywords=smoby+frisiertisch+gabby&qid=1729845795&sprefix=smoby+friseurtisch+gabby%2Caps%2C96&sr=8-3'
$ws.Cells.Item(35, 4).Value = '41 EUR'

$ws.Cells.Item(36, 1).Value = 'Laptop'
$ws.Cells.Item(36, 2).Value = 'https://m.media-amazon.com/images/I/91NEu43mBkL._AC_SX679_.jpg'
$ws.Cells.Item(36, 3).Value = 'https://www.amazon.de/-/en/VTech-Paw-Patrol-Learning-Laptop/dp/B09W2SYGY9/ref=sr_1_4?crid=3TXW07LSNODXD&dib=eyJ2IjoiMSJ9.aqSg65j_A-3LL_L4kXqS8Z_fzLKbDHn6lFls4301XS3w-4coU9jdO1i2ADl7ppvlygIJ6ehd_IEpZJq5Kx8EYaSHkTkRSqS2AMMpLVyCE-Cf9WQ-c_l8PgzTpcgN5gdREJmdy_b4kcrHXjsPXeI4pD40u3Wb5K-k-AdkzDvlePeRzqUkTZXrxx_EgJlHvf9_0n5BYIZKk1RVCtd6i3sSrfoyLFXsAr9Ca_VnK4aGH8wi8ylEztw9PZw5q7bnplYPkAjxRt1vJVt7nXqP8oFB0xKwI1GSctxexOhPVSPL6F4.hmsfZ-D44Jtm0oT6zTfoyJwbfBDYofU43mmZrk0f3kE&dib_tag=se&keywords=paw+patrol+lernlaptop&qid=1729845893&sprefix=paw+patrol+lern%2Caps%2C125&sr=8-4'
$ws.Cells.Item(36, 4).Value = '26 EUR'

$ws.Cells.Item(37, 1).Value = '(Clara) Animale Magice'
$ws.Cells.Item(37, 2).Value = 'https://m.media-amazon.com/images/I/81KzDuIaCNL._AC_SX679_.jpg'
$ws.Cells.Item(37, 3).Value = 'https://www.amazon.de/-/en/Bitzee-2024-Magicals-Interactive-Electronic/dp/B0CSND4V8T/ref=sr_1_1?crid=2FIGW3W6Z998O&dib=eyJ2IjoiMSJ9.i3mbxsNJzD4PPTehT-WNGjoVOceZeJn5R2mYk6ZG8YrB9dE4lVzmuJwKTd1OhvoCgrU31lzXEItydrB4bOfWY-b3Fwh3jyRkVAS62DTWucoMvk4nuAa9qLXJfLoDJeeUyVqzHTomwBHs0X0UAJl5f2pos3DC7uOdQOaysqO-Ykg9mkRHlldMMuQgxXs5-W7fADQCuX_rF_4jPbQd_oKD5LI4KMhPaCzbdXgx6pyQ_8c.5rf5mm466Ld8kUyhMQOTGNQQJRjQSwXkrFLOSdcnnvg&dib_tag=se&keywords=bitzee+magicals&qid=1729845929&sprefix=bitzee+ma%2Caps%2C140&sr=8-1'
$ws.Cells.Item(37, 4).Value = '30 EUR'

$ws.Range("C41").Select()